$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (FilesTab) query: change "    order by f.file_name" -> "     order By f.file_name ASC LIMIT 100"
$b4 = $ws.Range("B4").Value2
$b4 = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $b4

# Row 3 (SamplesTab) query: append order by clause
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`norder By samp.sample_id ASC LIMIT 100"

# Row 2 (CasesTab) query: append ORDER BY clause
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "

# Adjust wrapped row heights to reflect the longer text (auto height in real Excel)
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Update selection / active cell / scroll position to match the saved view state
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$null = $ws.Range("B4").Select()
